# Actualización de datos obtenidos el 6 de abril de 2016
# Updates the metadata sheet so that the age/sex/month-year fields are
# correctly classified as "dimension" (iaest-dimension: / dim) instead of
# "measure" (iaest-measure: / medida), switches the datatype for those
# columns from xsd:string to skos:Concept, and records the mapping files
# used for the age-group and sex code lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: semantic identifiers — these three were mis-tagged as "measure",
# they are actually SDMX/IAEST dimensions.
$ws.Range("A3").Value = "iaest-dimension:edad-grupos-quinquenales"
$ws.Range("G3").Value = "iaest-dimension:sexo"
$ws.Range("H3").Value = "iaest-dimension:mes-y-ano"

# Row 4: role flags follow the same three columns from "medida" to "dim".
$ws.Range("A4").Value = "dim"
$ws.Range("G4").Value = "dim"
$ws.Range("H4").Value = "dim"

# Row 5: datatype for the (now) dimension columns becomes a SKOS concept
# reference instead of a plain string. H5 (Mes y año datatype) is unchanged.
$ws.Range("A5").Value = "skos:Concept"
$ws.Range("G5").Value = "skos:Concept"

# Row 6 (new): mapping workbook for the age-group dimension, plus the same
# mapping reference for Sexo under column G.
$ws.Range("A6").Value = "mapping-edad-grupos-quinquenales.xlsx"
$ws.Range("G6").Value = "mapping-sexo.xlsx"

# New row 6 cells need the same style (Arial 10, style index 1) as the rest
# of the sheet — copy formatting from the row above.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)
